$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.27977745435881
$ws.Range("C2").Value = 20.01513248294323
$ws.Range("D2").Value = 5.999011257883991
$ws.Range("E2").Value = 16.3599649701172
$ws.Range("G2").Value = 3.658709265992564
$ws.Range("I2").Value = 43.39018794853718

$ws.Range("B3").Value = 16.80845620542125
$ws.Range("C3").Value = 19.03095298171618
$ws.Range("D3").Value = 5.886341013170409
$ws.Range("E3").Value = 15.42795899903451
$ws.Range("G3").Value = 3.665774595863357
$ws.Range("I3").Value = 41.74705059999795

$ws.Range("B4").Value = 16.52495127934529
$ws.Range("C4").Value = 18.40871444625651
$ws.Range("D4").Value = 5.818415080341074
$ws.Range("E4").Value = 14.83367852988898
$ws.Range("G4").Value = 3.670307741196599
$ws.Range("I4").Value = 40.70907039379494

$ws.Range("B5").Value = 16.41112971125974
$ws.Range("C5").Value = 18.15101172736924
$ws.Range("D5").Value = 5.79108615436758
$ws.Range("E5").Value = 14.58623719477642
$ws.Range("G5").Value = 3.672204432804039
$ws.Range("I5").Value = 40.27924060416047

$ws.Range("B6").Value = 16.39233964606281
$ws.Range("C6").Value = 18.1079827656509
$ws.Range("D6").Value = 5.786570477996603
$ws.Range("E6").Value = 14.54484060408023
$ws.Range("G6").Value = 3.672522371069423
$ws.Range("I6").Value = 40.20746904849166

$ws.Range("B7").Value = 16.52340901924836
$ws.Range("C7").Value = 18.40525518096005
$ws.Range("D7").Value = 5.818045042573145
$ws.Range("E7").Value = 14.83036237874036
$ws.Range("G7").Value = 3.670333120152682
$ws.Range("I7").Value = 40.70330061323982

$ws.Range("B8").Value = 17.11617998988317
$ws.Range("C8").Value = 19.67972215971064
$ws.Range("D8").Value = 5.959924059084231
$ws.Range("E8").Value = 16.04332735927769
$ws.Range("G8").Value = 3.661105138657576
$ws.Range("I8").Value = 42.82993721464702

$ws.Range("B9").Value = 18.31549519984824
$ws.Range("C9").Value = 22.02265412150252
$ws.Range("D9").Value = 6.246502329985754
$ws.Range("E9").Value = 18.28843694282577
$ws.Range("G9").Value = 3.644539644950052
$ws.Range("I9").Value = 46.75142658890895

$ws.Range("B10").Value = 19.2065254863081
$ws.Range("C10").Value = 23.63351757161559
$ws.Range("D10").Value = 6.460026389903472
$ws.Range("E10").Value = 19.95906145286519
$ws.Range("G10").Value = 3.633278346149001
$ws.Range("I10").Value = 49.46012747453848

$ws.Range("B11").Value = 19.6115856878747
$ws.Range("C11").Value = 24.33997384599312
$ws.Range("D11").Value = 6.557376378256579
$ws.Range("E11").Value = 20.67905300416601
$ws.Range("G11").Value = 3.628347535138037
$ws.Range("I11").Value = 50.65146846519302

$ws.Range("B12").Value = 19.76474781476391
$ws.Range("C12").Value = 24.60355145842846
$ws.Range("D12").Value = 6.594236952964063
$ws.Range("E12").Value = 20.94599803774185
$ws.Range("G12").Value = 3.626507553955342
$ws.Range("I12").Value = 51.09648501369985

$ws.Range("B13").Value = 19.73177443162579
$ws.Range("C13").Value = 24.54696278920632
$ws.Range("D13").Value = 6.586299074406537
$ws.Range("E13").Value = 20.8887593715783
$ws.Range("G13").Value = 3.626902623121028
$ws.Range("I13").Value = 51.00091833835749

$ws.Range("B14").Value = 19.62419195582505
$ws.Range("C14").Value = 24.36173832431566
$ws.Range("D14").Value = 6.560409189407872
$ws.Range("E14").Value = 20.70112868693547
$ws.Range("G14").Value = 3.628195615776998
$ws.Range("I14").Value = 50.68820412874621

$ws.Range("B15").Value = 19.55825997422021
$ws.Range("C15").Value = 24.24776551538222
$ws.Range("D15").Value = 6.54454938946932
$ws.Range("E15").Value = 20.58545862389973
$ws.Range("G15").Value = 3.628991142619525
$ws.Range("I15").Value = 50.49585426078749

$ws.Range("B16").Value = 19.18003208965117
$ws.Range("C16").Value = 23.5868051316743
$ws.Range("D16").Value = 6.453665679916858
$ws.Range("E16").Value = 19.91120795610494
$ws.Range("G16").Value = 3.633604415881942
$ws.Range("I16").Value = 49.38142673378726

$ws.Range("B17").Value = 18.94779421954244
$ws.Range("C17").Value = 23.17446250753899
$ws.Range("D17").Value = 6.39794165179753
$ws.Range("E17").Value = 19.48737127744995
$ws.Range("G17").Value = 3.63648340512393
$ws.Range("I17").Value = 48.68711171874063

$ws.Range("B18").Value = 18.81420153768831
$ws.Range("C18").Value = 22.9348225173079
$ws.Range("D18").Value = 6.365913070579499
$ws.Range("E18").Value = 19.23982517355283
$ws.Range("G18").Value = 3.63815742092728
$ws.Range("I18").Value = 48.28392615208452

$ws.Range("B19").Value = 18.76897258810743
$ws.Range("C19").Value = 22.85326529323308
$ws.Range("D19").Value = 6.355073681696302
$ws.Range("E19").Value = 19.15536192151856
$ws.Range("G19").Value = 3.638727334243913
$ws.Range("I19").Value = 48.14676396384805

$ws.Range("B20").Value = 18.97251934791591
$ws.Range("C20").Value = 23.21861408593599
$ws.Range("D20").Value = 6.403871508915744
$ws.Range("E20").Value = 19.53287887186566
$ws.Range("G20").Value = 3.636175061679776
$ws.Range("I20").Value = 48.761421408968

$ws.Range("B21").Value = 19.65579905944732
$ws.Range("C21").Value = 24.41625131778789
$ws.Range("D21").Value = 6.56801404840949
$ws.Range("E21").Value = 20.7563946713445
$ws.Range("G21").Value = 3.627815097252057
$ws.Range("I21").Value = 50.78022368024874

$ws.Range("B22").Value = 20.10097287074603
$ws.Range("C22").Value = 25.17594710745897
$ws.Range("D22").Value = 6.675253800925452
$ws.Range("E22").Value = 21.52283143072546
$ws.Range("G22").Value = 3.622509786031667
$ws.Range("I22").Value = 52.06387700725943

$ws.Range("B23").Value = 19.86355976320392
$ws.Range("C23").Value = 24.7726342471374
$ws.Range("D23").Value = 6.618032411755613
$ws.Range("E23").Value = 21.11679120042931
$ws.Range("G23").Value = 3.625326970503754
$ws.Range("I23").Value = 51.38210925785415

$ws.Range("B24").Value = 18.96134133992843
$ws.Range("C24").Value = 23.19866117739816
$ws.Range("D24").Value = 6.401190593222051
$ws.Range("E24").Value = 19.51231696738881
$ws.Range("G24").Value = 3.636314405017978
$ws.Range("I24").Value = 48.72783848384437

$ws.Range("B25").Value = 17.98843904134237
$ws.Range("C25").Value = 21.4073012330288
$ws.Range("D25").Value = 6.16830193886785
$ws.Range("E25").Value = 17.66496869973193
$ws.Range("G25").Value = 3.648859621270147
$ws.Range("I25").Value = 45.71951703758512
